{"js": "// Split three run-on \"numbered list\" / reference paragraphs into manual\n// line breaks (w:br) within a single run, one segment per line, matching\n// the author's re-flow of the \"Programa\" (PT/EN) and \"Bibliografia\"\n// paragraphs.\n//\n// Word (and Office.js) represents a manual line break typed with\n// Shift+Enter as <w:br/>. Feeding the U+000B \"vertical tab\" character\n// (the same code point Word itself uses internally for manual line\n// breaks, and what Range.Text/Paragraph.Text surface them as) to\n// insertText() makes the host create a single run containing\n// alternating <w:t>/<w:br/> children - exactly the structure produced\n// by pressing Shift+Enter between each segment.\n\nconst LINE_BREAK = \"\\u000b\";\n\nconst segments = {\n  pt: [\n    \"1. Introdu\u00e7\u00e3o aos sistemas de manufatura;\",\n    \"    2. Ind\u00fastria 4.0;\",\n    \"    3. Inteligencia artifical (ai), internet das coisas (IoT) e sistemas ciberf\u00edsicos;\",\n    \"    4. Sistemas de controle industrial; conceito de PLM e integra\u00e7\u00e3o com os sistemas de gest\u00e3o;\",\n    \"    5. Monitoramento e supervis\u00e3o de processos de produ\u00e7\u00e3o. Sistemas de controle da produ\u00e7\u00e3o, manufatura sustent\u00e1vel;\",\n    \"    6. Componentes de hardware para automa\u00e7\u00e3o de processos: controle num\u00e9rico, programa\u00e7\u00e3o CNC, controle discreto utilizando controladores l\u00f3gico program\u00e1veis e sistemas on-chip;\",\n    \"    7. Rob\u00f3tica industrial \u2013 programa\u00e7\u00e3o de rob\u00f4s e rob\u00f4s colaborativos;\",\n    \"    8. Sistemas de transporte de materiais e sistemas de armazenamento;\",\n    \"    9. Identifica\u00e7\u00e3o autom\u00e1tica e captura de dados \u2013 tecnologias de inspe\u00e7\u00e3o.\",\n  ],\n  en: [\n    \"1. Introduction to manufacturing systems;\",\n    \"2. 4.0 industry;\",\n    \"3. Artificial intelligence (AI), internet of things (IoT), and cyber-physical systems;\",\n    \"4. Industrial control systems; PLM concept and integration with management systems;\",\n    \"5. Production process monitoring and supervision. Production control systems, sustainable manufacturing;\",\n    \"6. Hardware components for process automation: numerical control, CNC programming, discrete control using programmable logic controllers and on-chip systems;\",\n    \"7. Industrial robotics \u2013 programming of robots and collaborative robots;\",\n    \"8. Material transportation systems and storage systems;\",\n    \"9. Automatic identification and data capture \u2013 inspection technologies.\",\n  ],\n  bib: [\n    \"GROOVER, M.P. Automa\u00e7\u00e3o Industrial e Sistemas de Manufatura, 561p., 3a Edi\u00e7\u00e3o - S\u00e3o Paulo, Pearson Prentice Hall, 2011.\",\n    \"RICHARD, L.S., ERNEST, L.H., Handbook of Industrial Automation, Marcel Dekker,Inc. NewYork, 2000.\",\n    \"ADALBERTO FILHO ET. AL, Automa\u00e7\u00e3o & Sociedade: Quarta revolu\u00e7\u00e3o Industrial, um olhar para o Brasil, 1a Edi\u00e7\u00e3o, Brasport Livros e Multim\u00eddia Limitada.\",\n  ],\n};\n\nconst targets = [\n  { key: \"pt\", prefix: \"1. Introdu\u00e7\u00e3o aos sistemas de manufatura;\" },\n  { key: \"en\", prefix: \"1. Introduction to manufacturing systems;\" },\n  { key: \"bib\", prefix: \"GROOVER, M.P. Automa\u00e7\u00e3o Industrial e Sistemas de Manufatura\" },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of targets) {\n  const paragraph = paragraphs.items.find((p) => p.text.indexOf(target.prefix) === 0);\n  if (!paragraph) {\n    continue;\n  }\n\n  // Preserve the paragraph's run-level formatting (e.g. the italic\n  // English translation) before wiping its contents.\n  const probeFont = paragraph.getRange().font;\n  probeFont.load(\"italic,bold,underline\");\n  await context.sync();\n  const wasItalic = probeFont.italic;\n  const wasBold = probeFont.bold;\n  const wasUnderline = probeFont.underline;\n\n  paragraph.clear();\n  await context.sync();\n\n  const newText = segments[target.key].join(LINE_BREAK);\n  const insertedRange = paragraph.insertText(newText, Word.InsertLocation.start);\n  if (wasItalic) {\n    insertedRange.font.italic = true;\n  }\n  if (wasBold) {\n    insertedRange.font.bold = true;\n  }\n  if (wasUnderline && wasUnderline !== \"None\") {\n    insertedRange.font.underline = wasUnderline;\n  }\n  await context.sync();\n}\n", "ps1": "# Split three run-on \"numbered list\" / reference paragraphs into manual\n# line breaks, one segment per line, matching the author's re-flow of the\n# \"Programa\" (PT/EN) and \"Bibliografia\" paragraphs.\n#\n# [char]11 is the \"vertical tab\" code point - the same one Word itself\n# uses for a manual line break (Shift+Enter). Assigning a Range's .Text\n# with embedded [char]11 characters makes Word split the run into\n# alternating text/<w:br/> runs of content while keeping the existing\n# run formatting (e.g. italics) intact, exactly like typing Shift+Enter\n# between each segment.\n\n$d = $word.ActiveDocument\n\n$lineBreak = [char]11\n\n$ptSegments = @(\n    \"1. Introdu\u00e7\u00e3o aos sistemas de manufatura;\",\n    \"    2. Ind\u00fastria 4.0;\",\n    \"    3. Inteligencia artifical (ai), internet das coisas (IoT) e sistemas ciberf\u00edsicos;\",\n    \"    4. Sistemas de controle industrial; conceito de PLM e integra\u00e7\u00e3o com os sistemas de gest\u00e3o;\",\n    \"    5. Monitoramento e supervis\u00e3o de processos de produ\u00e7\u00e3o. Sistemas de controle da produ\u00e7\u00e3o, manufatura sustent\u00e1vel;\",\n    \"    6. Componentes de hardware para automa\u00e7\u00e3o de processos: controle num\u00e9rico, programa\u00e7\u00e3o CNC, controle discreto utilizando controladores l\u00f3gico program\u00e1veis e sistemas on-chip;\",\n    \"    7. Rob\u00f3tica industrial \u2013 programa\u00e7\u00e3o de rob\u00f4s e rob\u00f4s colaborativos;\",\n    \"    8. Sistemas de transporte de materiais e sistemas de armazenamento;\",\n    \"    9. Identifica\u00e7\u00e3o autom\u00e1tica e captura de dados \u2013 tecnologias de inspe\u00e7\u00e3o.\"\n)\n\n$enSegments = @(\n    \"1. Introduction to manufacturing systems;\",\n    \"2. 4.0 industry;\",\n    \"3. Artificial intelligence (AI), internet of things (IoT), and cyber-physical systems;\",\n    \"4. Industrial control systems; PLM concept and integration with management systems;\",\n    \"5. Production process monitoring and supervision. Production control systems, sustainable manufacturing;\",\n    \"6. Hardware components for process automation: numerical control, CNC programming, discrete control using programmable logic controllers and on-chip systems;\",\n    \"7. Industrial robotics \u2013 programming of robots and collaborative robots;\",\n    \"8. Material transportation systems and storage systems;\",\n    \"9. Automatic identification and data capture \u2013 inspection technologies.\"\n)\n\n$bibSegments = @(\n    \"GROOVER, M.P. Automa\u00e7\u00e3o Industrial e Sistemas de Manufatura, 561p., 3a Edi\u00e7\u00e3o - S\u00e3o Paulo, Pearson Prentice Hall, 2011.\",\n    \"RICHARD, L.S., ERNEST, L.H., Handbook of Industrial Automation, Marcel Dekker,Inc. NewYork, 2000.\",\n    \"ADALBERTO FILHO ET. AL, Automa\u00e7\u00e3o & Sociedade: Quarta revolu\u00e7\u00e3o Industrial, um olhar para o Brasil, 1a Edi\u00e7\u00e3o, Brasport Livros e Multim\u00eddia Limitada.\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    if ($t.StartsWith(\"1. Introdu\") -and $t.Contains(\"manufatura\")) {\n        $p.Range.Text = [string]::Join($lineBreak, $ptSegments)\n    }\n    elseif ($t.StartsWith(\"1. Introduction to manufacturing\")) {\n        $p.Range.Text = [string]::Join($lineBreak, $enSegments)\n    }\n    elseif ($t.StartsWith(\"GROOVER\")) {\n        $p.Range.Text = [string]::Join($lineBreak, $bibSegments)\n    }\n}\n"}
